# Auto-generated update of H:N numeric columns per the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 12998.6045
$ws.Range("I132").Value = 1768.125
$ws.Range("J132").Value = 102842.445
$ws.Range("K132").Value = 5304.375
$ws.Range("L132").Value = 308527.335
$ws.Range("M132").Value = -2774.375
$ws.Range("N132").Value = -313587.335

$ws.Range("H135").Value = 8621687
$ws.Range("I135").Value = 674.561
$ws.Range("J135").Value = 29413542
$ws.Range("K135").Value = 6071.049
$ws.Range("L135").Value = 264721878
$ws.Range("M135").Value = -3536.049
$ws.Range("N135").Value = -264726948

$ws.Range("H136").Value = 36524.25
$ws.Range("J136").Value = 36524.25
$ws.Range("L136").Value = 36524.25
$ws.Range("N136").Value = -46724.25

$ws.Range("H137").Value = 2881.303
$ws.Range("I137").Value = 875.7907
$ws.Range("J137").Value = 6630.7393
$ws.Range("K137").Value = 2627.3721
$ws.Range("L137").Value = 19892.2179
$ws.Range("M137").Value = -77.37210000000005
$ws.Range("N137").Value = -24992.2179

$ws.Range("H138").Value = 2636.5386
$ws.Range("I138").Value = 1412
$ws.Range("J138").Value = 4595.8
$ws.Range("K138").Value = 4236
$ws.Range("L138").Value = 13787.4
$ws.Range("M138").Value = 904
$ws.Range("N138").Value = -24067.4

$ws.Range("H141").Value = 912.7857
$ws.Range("J141").Value = 6170
$ws.Range("L141").Value = 18510
$ws.Range("N141").Value = -28870

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1538.9508
$ws.Range("I61").Value = 1354.2
$ws.Range("J61").Value = 2058.5625
$ws.Range("K61").Value = 1354.2
$ws.Range("L61").Value = 2058.5625
$ws.Range("M61").Value = -1142.2
$ws.Range("N61").Value = -2482.5625

$ws.Range("H122").Value = 1548.2115
$ws.Range("I122").Value = 1590.439
$ws.Range("K122").Value = 4771.317
$ws.Range("M122").Value = -2321.317

$ws.Range("H132").Value = 8773901
$ws.Range("I132").Value = 15626326
$ws.Range("J132").Value = 2795.56
$ws.Range("K132").Value = 46878978
$ws.Range("L132").Value = 8386.68
$ws.Range("M132").Value = -46876448
$ws.Range("N132").Value = -13446.68

$ws.Range("H136").Value = 1538.9508
$ws.Range("I136").Value = 1354.2
$ws.Range("J136").Value = 2058.5625
$ws.Range("K136").Value = 4062.6
$ws.Range("L136").Value = 6175.6875
$ws.Range("M136").Value = -1512.6
$ws.Range("N136").Value = -11275.6875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 42800
$ws.Range("J95").Value = 42800
$ws.Range("L95").Value = 42800
$ws.Range("N95").Value = -48292

$ws.Range("H99").Value = 2526.8948
$ws.Range("I99").Value = 2200
$ws.Range("J99").Value = 3752.75
$ws.Range("K99").Value = 2200
$ws.Range("L99").Value = 3752.75
$ws.Range("M99").Value = -702
$ws.Range("N99").Value = -6748.75

$ws.Range("H134").Value = 2466.9697
$ws.Range("I134").Value = 869.36584
$ws.Range("J134").Value = 3596.3103
$ws.Range("K134").Value = 2608.09752
$ws.Range("L134").Value = 10788.9309
$ws.Range("M134").Value = -73.09752000000026
$ws.Range("N134").Value = -15858.9309

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2286.73
$ws.Range("I31").Value = 973.617
$ws.Range("J31").Value = 3451.1887
$ws.Range("K31").Value = 973.617
$ws.Range("L31").Value = 3451.1887
$ws.Range("M31").Value = -678.617
$ws.Range("N31").Value = -4041.1887

$ws.Range("H34").Value = 2286.73
$ws.Range("I34").Value = 973.617
$ws.Range("J34").Value = 3451.1887
$ws.Range("K34").Value = 973.617
$ws.Range("L34").Value = 3451.1887
$ws.Range("M34").Value = -771.617
$ws.Range("N34").Value = -3855.1887

$ws.Range("H58").Value = 1975.16
$ws.Range("I58").Value = 1378.1052
$ws.Range("K58").Value = 1378.1052
$ws.Range("M58").Value = -1175.1052

$ws.Range("H132").Value = 37826.2
$ws.Range("I132").Value = 2290.5
$ws.Range("J132").Value = 144433.3
$ws.Range("K132").Value = 6871.5
$ws.Range("L132").Value = 433299.9
$ws.Range("M132").Value = -4341.5
$ws.Range("N132").Value = -438359.9

$ws.Range("H134").Value = 312711.03
$ws.Range("I134").Value = 1052.1666
$ws.Range("J134").Value = 936028.75
$ws.Range("K134").Value = 3156.4998
$ws.Range("L134").Value = 2808086.25
$ws.Range("M134").Value = -621.4998000000001
$ws.Range("N134").Value = -2813156.25

$ws.Range("H136").Value = 1975.16
$ws.Range("I136").Value = 1378.1052
$ws.Range("K136").Value = 4134.3156
$ws.Range("M136").Value = -1584.3156

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4042.0344
$ws.Range("I113").Value = 5804.684
$ws.Range("K113").Value = 17414.052
$ws.Range("M113").Value = -15244.052

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 926830.6
$ws.Range("I2").Value = 1506070
$ws.Range("J2").Value = 47.6
$ws.Range("K2").Value = 1506070
$ws.Range("L2").Value = 47.6
$ws.Range("M2").Value = -1505957
$ws.Range("N2").Value = -273.6

$ws.Range("H113").Value = 1560.5714
$ws.Range("I113").Value = 1502.4445
$ws.Range("K113").Value = 1502.4445
$ws.Range("M113").Value = 667.5554999999999

$ws.Range("H123").Value = 12930.637
$ws.Range("J123").Value = 12930.637
$ws.Range("L123").Value = 12930.637
$ws.Range("N123").Value = -17830.637

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3175.1738
$ws.Range("I16").Value = 3101.875
$ws.Range("J16").Value = 3342.7144
$ws.Range("K16").Value = 3101.875
$ws.Range("L16").Value = 3342.7144
$ws.Range("M16").Value = -2931.875
$ws.Range("N16").Value = -3682.7144

$ws.Range("H93").Value = 1624.0869
$ws.Range("I93").Value = 1475.8334
$ws.Range("J93").Value = 1676.4117
$ws.Range("K93").Value = 1475.8334
$ws.Range("L93").Value = 1676.4117
$ws.Range("M93").Value = -227.8334
$ws.Range("N93").Value = -4172.411700000001

$ws.Range("H122").Value = 2027.0952
$ws.Range("I122").Value = 1840.6428
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 5521.928400000001
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -3071.928400000001
$ws.Range("N122").Value = -12100

$ws.Range("H132").Value = 2336
$ws.Range("I132").Value = 1440
$ws.Range("K132").Value = 4320
$ws.Range("M132").Value = -1790

$ws.Range("H136").Value = 1441.4637
$ws.Range("I136").Value = 1121.7407
$ws.Range("J136").Value = 2592.4666
$ws.Range("K136").Value = 3365.2221
$ws.Range("L136").Value = 7777.399800000001
$ws.Range("M136").Value = -815.2221
$ws.Range("N136").Value = -12877.3998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 33334364

$ws.Range("H113").Value = 868.61536
$ws.Range("I113").Value = 789.2
$ws.Range("J113").Value = 1133.3334
$ws.Range("K113").Value = 2367.6
$ws.Range("L113").Value = 3400.0002
$ws.Range("M113").Value = -197.6000000000004
$ws.Range("N113").Value = -7740.0002

$ws.Range("H132").Value = 1555.2273
$ws.Range("I132").Value = 1503.4822
$ws.Range("J132").Value = 1645.7812
$ws.Range("K132").Value = 4510.446599999999
$ws.Range("L132").Value = 4937.3436
$ws.Range("M132").Value = -1980.446599999999
$ws.Range("N132").Value = -9997.3436

$ws.Range("H136").Value = 689.0513
$ws.Range("I136").Value = 486.94232
$ws.Range("J136").Value = 1093.2693
$ws.Range("K136").Value = 1460.82696
$ws.Range("L136").Value = 3279.8079
$ws.Range("M136").Value = 1089.17304
$ws.Range("N136").Value = -8379.8079
